$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Header row (row 1): the three date columns G1 (Vehicle_Date of
# Manufacture), M1 (Insurant_birthdate) and T1 (Product_startdate) get a
# Text number format applied so the whole column can safely hold literal
# date strings; labels themselves are unchanged. ---
$ws.Range("G1").NumberFormat = "@"
$ws.Range("M1").NumberFormat = "@"
$ws.Range("T1").NumberFormat = "@"

# --- Apply the same Text format to the data cells in those columns for
# rows 2-4 ---
$ws.Range("G2").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("M2").NumberFormat = "@"
$ws.Range("M3").NumberFormat = "@"
$ws.Range("M4").NumberFormat = "@"
$ws.Range("T2").NumberFormat = "@"
$ws.Range("T3").NumberFormat = "@"
$ws.Range("T4").NumberFormat = "@"

# --- Replace the real Excel date values with literal date text values
# (data-driven test case completed: dates + insurance sums filled in) ---
$ws.Range("G4").Value = "04/09/1989"
$ws.Range("G3").Value = "04/09/2001"
$ws.Range("G2").Value = "04/09/2019"

$ws.Range("M2").Value = "12/12/1989"
$ws.Range("M3").Value = "12/12/1977"
$ws.Range("M4").Value = "12/12/1956"

$ws.Range("T2").Value = "15/06/2023"
$ws.Range("T3").Value = "15/06/2023"
$ws.Range("T4").Value = "15/06/2023"

# --- Product_insurancesum (U column): was free-text European formatted
# numbers, now a plain numeric value for every row ---
$ws.Range("U2").Value = 3000000
$ws.Range("U3").Value = 3000000
$ws.Range("U4").Value = 3000000

# --- A handful of stray formatted-but-empty cells appear below the table
# (rows 7-11), reusing the pre-existing dd/mm/yyyy date style. These are
# left-over artifacts from entering/removing extra test rows. ---
$ws.Range("S7").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("L8").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("S8").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("F9").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("L9").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("S9").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("F10").NumberFormat = "dd\/mm\/yyyy"
$ws.Range("L10").NumberFormat = "dd\/mm\/yyyy"

$ws.Range("F11").NumberFormat = "dd\/mm\/yyyy"
